$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(2, 3, 5, 6, 7, 12, 13, 14, 15)  # B, C, E, F, G, L, M, N, O

$data = @{
    2 = @(0.9646773129393296, 0.2335068093486257, 0.1166771519514311, 0.4443680307746263, 0.002438830678968888, 0.1888783094692243, 0.2188512100700279, 1.397437926242375, 2.820092277855082)
    3 = @(0.8867854317317949, 0.2248336650023077, 0.1176737342298653, 0.387822817061874, 0.002441577148969299, 0.1862223997423698, 0.2063987760370196, 1.410139803999641, 2.831181935501405)
    4 = @(0.8391672877897918, 0.219466063207534, 0.118325077472484, 0.3531389305168915, 0.002443354212374883, 0.1846825474385696, 0.1988285551893441, 1.418488560385882, 2.83987984088958)
    5 = @(0.8198158797700899, 0.2172682230444565, 0.1186004339467908, 0.3390132514313251, 0.002444101260552418, 0.1840779554454954, 0.1957628285134589, 1.422029092433178, 2.843898766760361)
    6 = @(0.8166058498846667, 0.2169026429942136, 0.1186467567051839, 0.336668177824194, 0.002444226691133695, 0.1839789487872281, 0.1952549319681793, 1.422625354046126, 2.844594747939453)
    7 = @(0.8389060903420784, 0.2194364647055949, 0.1183287508087973, 0.3529483938344953, 0.002443364194836452, 0.1846743008620848, 0.1987871317445169, 1.418535748939171, 2.83993212112847)
    8 = @(0.9377777141388037, 0.2305251529322163, 0.1170125954339136, 0.4248636149813336, 0.002439758873974304, 0.1879437205041441, 0.2145420208262507, 1.401703496057756, 2.823523663841684)
    9 = @(1.133274078586567, 0.251930146560511, 0.1147439889502682, 0.5661985755041457, 0.002433405563133223, 0.1950742716254013, 0.246030998258064, 1.373052834495418, 2.806358946201811)
    10 = @(1.277848428145433, 0.2674444164516956, 0.1132668919110014, 0.6702781546542269, 0.002429170355258515, 0.2007497736970123, 0.2695218570393507, 1.354653180082408, 2.802939134414288)
    11 = @(1.343816721873566, 0.2744553074152805, 0.1126359362861873, 0.7176906081379002, 0.002427336650809451, 0.2034262046748552, 0.2802846205498071, 1.346856860204461, 2.803387872789699)
    12 = @(1.368825126872309, 0.2771033345299259, 0.1124028911614854, 0.7356546913071611, 0.002426655565689257, 0.2044532610834864, 0.2843710716614893, 1.343987024804377, 2.803846672986367)
    13 = @(1.363437909266622, 0.2765333411764459, 0.1124528200497297, 0.7317853510981394, 0.002426801658928333, 0.2042314643520768, 0.2834905029294674, 1.344601428390185, 2.803735004306816)
    14 = @(1.345873631336474, 0.2746733001869472, 0.1126166456506148, 0.7191683204515869, 0.00242728035146049, 0.2035104299787207, 0.2806205994949309, 1.346619104713618, 2.803419825465312)
    15 = @(1.335118567941947, 0.273533075786645, 0.1127177595683477, 0.7114413442032514, 0.002427575294493778, 0.2030705385749059, 0.2788641048298999, 1.347865727369467, 2.803264408010079)
    16 = @(1.273541179529786, 0.2669852879763823, 0.1133089502392135, 0.6671810134426437, 0.002429292056483663, 0.2005767618066443, 0.2688200107935401, 1.355174229674901, 2.802950197114995)
    17 = @(1.235815983464761, 0.2629563898316292, 0.1136821167447559, 0.6400460337125793, 0.002430368986956316, 0.1990711044923614, 0.2626777773741722, 1.359804682265533, 2.803271275532779)
    18 = @(1.214136411564198, 0.2606346936558452, 0.1139006103342819, 0.6244449056556647, 0.002430997158295999, 0.1982139993191083, 0.2591521539384658, 1.3625220050492, 2.803644579461888)
    19 = @(1.206799375051162, 0.2598478591373521, 0.1139752514225503, 0.619163680173358, 0.002431211350879032, 0.1979253300763446, 0.2579596856226942, 1.363451321428919, 2.803803351360216)
    20 = @(1.239829941196206, 0.2633857273017952, 0.1136419933032196, 0.642933953830422, 0.002430253440962505, 0.1992304626299983, 0.2633308820888587, 1.359306173309896, 2.803217569252354)
    21 = @(1.351031944962244, 0.2752198262032266, 0.1125683664818391, 0.7228739723491628, 0.002427139387628537, 0.2037218479505611, 0.2814632674477195, 1.346024226628728, 2.803504556125148)
    22 = @(1.423869449763401, 0.2829141305697647, 0.1119009819407486, 0.7751780083420101, 0.002425181659919049, 0.2067361823445566, 0.2933768415748048, 1.337824338327344, 2.805376185495902)
    23 = @(1.384980411943559, 0.2788112383298937, 0.1122540428938752, 0.7472568307830727, 0.002426219466156522, 0.2051201693937372, 0.2870126421134387, 1.342156809414824, 2.80422295875286)
    24 = @(1.238015202324277, 0.2631916407327424, 0.1136601207906116, 0.6416283278902171, 0.002430305651385992, 0.1991583902819798, 0.2630355959008597, 1.359531377214225, 2.803241262079268)
    25 = @(1.080218769427063, 0.2461764113660934, 0.115324346733944, 0.5279251897347166, 0.002435048025260371, 0.1930684540135132, 0.23744946110142, 1.380337790160809, 2.809391237680074)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $ws.Cells.Item([int]$r, $cols[$i]).Value = $vals[$i]
    }
}
